$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 628, shifting existing rows (628..669) down to (629..670)
$ws.Rows.Item(628).Insert()

# Fill the newly inserted row with the new data point.
$dateCell = $ws.Cells.Item(628, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/11"
$dateCell.Style = "Normal"

$ws.Cells.Item(628, 2).Value = "日"
$ws.Cells.Item(628, 3).Value = 16
$ws.Cells.Item(628, 4).Value = 139
